# RPA datasets push 2023-12-05
# A new IPO entry ("IBKS스팩23호") is inserted as the latest record (row 2),
# pushing the existing rows down by one and dropping the oldest record
# (the row that used to be row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A="IBKS스팩23호";            B="2024.01.08~01.12"; C="2,400~2,800";    D="-";     E=16000;  F="NH투자증권"},
    @{Row=3;  A="하나스팩30호";             B="2024.01.05~01.11"; C="13,000~15,000";  D="-";     E=19500;  F="하나증권"},
    @{Row=4;  A="디에스단석(구,단석산업)";  B="2023.12.08~12.11"; C="2,000~2,000";    D="-";     E=14000;  F="하나증권"},
    @{Row=5;  A="블루엠텍";                 B="2023.12.05~12.11"; C="79,000~89,000";  D="-";     E=96380;  F="KB증권,NH투자증권"},
    @{Row=6;  A="LS머트리얼즈";             B="2023.11.22~11.28"; C="15,000~19,000";  D="19000"; E=21000;  F="하나증권,키움증권"},
    @{Row=7;  A="삼성스팩9호";              B="2023.11.22~11.28"; C="4,400~5,500";    D="6000";  E=64350;  F="키움증권,KB증권,이베스트투자증권,하이투자증권,NH투자증권"},
    @{Row=8;  A="교보스팩15호";             B="2023.11.20~11.21"; C="2,000~2,000";    D="2000";  E=20000;  F="삼성증권"},
    @{Row=9;  A="케이엔에스";               B="2023.11.20~11.21"; C="2,000~2,000";    D="2000";  E=7000;   F="교보증권"},
    @{Row=10; A="NH스팩30호";               B="2023.11.16~11.22"; C="19,000~22,000";  D="23000"; E=14250;  F="신영증권"},
    @{Row=11; A="와이바이오로직스";         B="2023.11.15~11.16"; C="2,000~2,000";    D="2000";  E=16000;  F="NH투자증권"},
    @{Row=12; A="에이텀";                   B="2023.11.10~11.16"; C="9,000~11,000";   D="9000";  E=13500;  F="유안타증권"},
    @{Row=13; A="에이에스텍";               B="2023.11.09~11.15"; C="23,000~30,000";  D="18000"; E=14950;  F="하나증권"},
    @{Row=14; A="그린리소스";               B="2023.11.07~11.13"; C="21,000~25,000";  D="28000"; E=29547;  F="미래에셋증권"},
    @{Row=15; A="한선엔지니어링";           B="2023.11.03~11.09"; C="11,000~14,000";  D="17000"; E=18040;  F="NH투자증권"},
    @{Row=16; A="에코아이";                 B="2023.11.02~11.08"; C="5,200~6,000";    D="7000";  E=22100;  F="대신증권"},
    @{Row=17; A="동인기연(유가)";           B="2023.11.01~11.07"; C="28,500~34,700";  D="34700"; E=59251;  F="KB증권"},
    @{Row=18; A="스톰테크";                 B="2023.11.01~11.07"; C="33,000~37,000";  D="30000"; E=60654;  F="NH투자증권"},
    @{Row=19; A="에코프로머티리얼즈(유가)"; B="2023.10.31~11.06"; C="8,000~9,500";    D="11000"; E=26800;  F="하이투자증권"},
    @{Row=20; A="캡스톤파트너스";           B="2023.10.30~11.03"; C="36,200~44,000";  D="36200"; E=524031; F="미래에셋증권,NH투자증권,하이투자증권"},
    @{Row=21; A="한국스팩13호";             B="2023.10.26~11.01"; C="3,200~3,600";    D="4000";  E=6384;   F="NH투자증권"}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}

Write-Output "Applied RPA datasets push 2023-12-05: updated rows 2-21"
